# AdjustableParameters.xlsx - "Add files via upload" (IRIS Imager v1.2 update)
#
# Adds two new Hydric-Soils-Technical-Standard parameter rows (Avg_Thickness /
# Avg_Depth), renames the Fe output file, adds a NOTES column, and adds
# data-validation rules to every adjustable cell on the Parameters sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# ------------------------------------------------------------------
# 1. Insert two new rows (18 & 19) for Avg_Thickness / Avg_Depth, pushing
#    the existing OutputRemoval..Image_Path rows down from 18-25 to 20-27.
# ------------------------------------------------------------------
$ws.Rows("18:19").Insert()

$ws.Range("A18").Value = "Avg_Thickness"
$ws.Range("B18").Value = 6
$ws.Range("B18").Style = "Input"
$ws.Range("C18").Value = "Hydric Soils Technical Standard: Thickness of zone to average removal across (in)"
$ws.Range("D18").Value = "0<x<Avg_Depth"

$ws.Range("A19").Value = "Avg_Depth"
$ws.Range("B19").Value = 12
$ws.Range("B19").Style = "Input"
$ws.Range("C19").Value = "Hydric Soils Technical Standard: Bottom of zone to average removal across (in)"
$ws.Range("D19").Value = "Avg_Thickness<x<length of film"

# ------------------------------------------------------------------
# 2. Rename the Fe output file (now on row 23, was row 21).
# ------------------------------------------------------------------
$ws.Range("B23").Value = "OutData_Fe.xlsx"

# ------------------------------------------------------------------
# 3. Add the NOTES column (F1:F7) describing the v1.2 workbook.
# ------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 48.59

# F1 reuses the same bold+underline header look as row 1 (copy format from D1).
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "NOTES"

$ws.Range("F2").Value = "This is the Adjustable Parameters file for IRIS Imager v1.2"
$ws.Range("F3").Value = "Adjust any values in column B"
$ws.Range("F4").Value = "Do not change tab names"
$ws.Range("F5").Value = "Do not add rows"
$ws.Range("F6").Value = "Save different versions of the file for different file types/analyses"
$ws.Range("F7").Value = "For more info, see the User Guide"

# ------------------------------------------------------------------
# 4. Data validation rules for every adjustable (column B) cell.
# ------------------------------------------------------------------

# DPI (B2): must be > 0
$val = $ws.Range("B2").Validation
$val.Add(7, 1, 1, "B2>0")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "DPI must be >0"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# x_Crop, y_Crop, CircleRad_min (B3:B4, B10): must be > 0
$val = $ws.Range("B3:B4").Validation
$val.Add(7, 1, 1, "B3>0")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be >0"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

$val = $ws.Range("B10").Validation
$val.Add(7, 1, 1, "B3>0")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be >0"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# Scale_Factor (B5): must be >0 and <=1
$val = $ws.Range("B5").Validation
$val.Add(7, 1, 1, "AND(B5>0,B5<=1)")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be >0 and <=1"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# Start_Upside_down_flag (B6): 0/1 list
$val = $ws.Range("B6").Validation
$val.Add(3, 1, 1, '"0,1"')
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# End_Upside_down_flag (B7): 0/1 list
$val = $ws.Range("B7").Validation
$val.Add(3, 1, 1, '"0,1"')
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# White_L, Fe_b_Threshold (B8:B9): between 0 and 100
$val = $ws.Range("B8:B9").Validation
$val.Add(7, 1, 1, "AND(B8>=0,B8<=100)")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be between 0 and 100"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# CircleRad_max (B11): must be >0 and > CircleRad_min
$val = $ws.Range("B11").Validation
$val.Add(7, 1, 1, "AND(B11>0,B11>B10)")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be >0 and > CircleRad_min"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# Crop_Buffer, Circle_NaN_buffer (B12:B13): must be >1
$val = $ws.Range("B12:B13").Validation
$val.Add(7, 1, 1, "B12>1")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be >1"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# n_circles_start, n_circles_end (B14:B15): whole number
$val = $ws.Range("B14:B15").Validation
$val.Add(1, 1, 1, 0, 1000000)
$val.ErrorTitle = "Out of range"
$val.ErrorMessage = "Must be integer"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# n_clusters_start, n_clusters_end (B16:B17): positive whole number
$val = $ws.Range("B16:B17").Validation
$val.Add(1, 1, 1, 1, 1000000)
$val.ErrorTitle = "Out of range"
$val.ErrorMessage = "Must be positive integer"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# Avg_Thickness (B18): must be between 0 and Avg_Depth
$val = $ws.Range("B18").Validation
$val.Add(7, 1, 1, "AND(B18>0,B18<B19)")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be between 0 and Avg_Depth"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# Avg_Depth (B19): must be > Avg_Thickness
$val = $ws.Range("B19").Validation
$val.Add(7, 1, 1, "B19>B18")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be >Avg_Thickness"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# OutputRemoval (B20): 1/0 list
$val = $ws.Range("B20").Validation
$val.Add(3, 1, 1, '"1,0"')
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# OutputRemovalDepth (B21): must be >=0
$val = $ws.Range("B21").Validation
$val.Add(7, 1, 1, "B21>=0")
$val.ErrorTitle = "Out of Range"
$val.ErrorMessage = "Must be >=0"
$val.IgnoreBlank = $true
$val.ShowInput = $true
$val.ShowError = $true

# IRIS_Type (B22): list, blank not allowed
$val = $ws.Range("B22").Validation
$val.Add(3, 1, 1, '"Fe, Mn, S"')
$val.IgnoreBlank = $false
$val.ShowInput = $true
$val.ShowError = $true

# Start_circle_color (B25): list, blank not allowed
$val = $ws.Range("B25").Validation
$val.Add(3, 1, 1, '"bright, dark"')
$val.IgnoreBlank = $false
$val.ShowInput = $true
$val.ShowError = $true

# End_circle_color (B26): list, blank not allowed
$val = $ws.Range("B26").Validation
$val.Add(3, 1, 1, '"dark, bright"')
$val.IgnoreBlank = $false
$val.ShowInput = $true
$val.ShowError = $true

# ------------------------------------------------------------------
# 5. Selection parks on B20 (matches the saved workbook state).
# ------------------------------------------------------------------
$ws.Range("B20").Select()
